$wb = $excel.ActiveWorkbook

# --- Status text: "Ready for handoff" -> "In Translation" -----------------
# "Overview" sheet holds per-locale status in columns E (zh-cn) / F (de-de)
# for the two content rows.
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "In Translation"
$wsOverview.Range("F2").Value = "In Translation"
$wsOverview.Range("E3").Value = "In Translation"
$wsOverview.Range("F3").Value = "In Translation"

# Each locale sheet repeats the same Status column (column C) in its table.
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "In Translation"
$wsZhCn.Range("C3").Value = "In Translation"

$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "In Translation"
$wsDeDe.Range("C3").Value = "In Translation"

# --- Column widths follow the shorter "In Translation" label --------------
# The report generator re-fits the Status columns to the new text, shrinking
# them from the old "Ready for handoff" width down to the new, narrower one.
$wsOverview.Range("E:E").ColumnWidth = 12.5
$wsOverview.Range("F:F").ColumnWidth = 12.5
$wsZhCn.Range("C:C").ColumnWidth = 12.5
$wsDeDe.Range("C:C").ColumnWidth = 12.5
